# Adds two new match rows (146 and 147) to the Ekstraklasa 2023-2024 sheet,
# mirroring the formatting of the last existing data row (145).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 146 -------------------------------------------------------------
$ws.Range("A145").Copy()
$ws.Range("A146").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A146").Value = 145

$ws.Range("B146").Value = "poland"
$ws.Range("C146").Value = "ekstraklasa"
$ws.Range("D146").Value = "2023-2024"

$ws.Range("E145").Copy()
$ws.Range("E146").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E146").Value = 45262.72916666666

$ws.Range("F146").Value = "Widzew Lodz"
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = "Radomiak Radom"
$ws.Range("I146").Value = 3
$ws.Range("J146").Value = 2.4
$ws.Range("K146").Value = "27/11/2023 19:12"
$ws.Range("L146").Value = 2.2
$ws.Range("M146").Value = "02/12/2023 17:22"
$ws.Range("N146").Value = 3.33
$ws.Range("O146").Value = "27/11/2023 19:12"
$ws.Range("P146").Value = 3.39
$ws.Range("Q146").Value = "02/12/2023 17:22"
$ws.Range("R146").Value = 3.08
$ws.Range("S146").Value = "27/11/2023 19:12"
$ws.Range("T146").Value = 3.52
$ws.Range("U146").Value = "02/12/2023 17:22"
$ws.Range("V146").Value = "https://www.betexplorer.com/football/poland/ekstraklasa/widzew-lodz-radomiak-radom/jwkespr7/"

# --- Row 147 -------------------------------------------------------------
$ws.Range("A145").Copy()
$ws.Range("A147").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A147").Value = 146

$ws.Range("B147").Value = "poland"
$ws.Range("C147").Value = "ekstraklasa"
$ws.Range("D147").Value = "2023-2024"

$ws.Range("E145").Copy()
$ws.Range("E147").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E147").Value = 45262.83333333334

$ws.Range("F147").Value = "Korona Kielce"
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = "Lech Poznan"
$ws.Range("I147").Value = 1
$ws.Range("J147").Value = 3.95
$ws.Range("K147").Value = "26/11/2023 18:13"
$ws.Range("L147").Value = 3.24
$ws.Range("M147").Value = "02/12/2023 19:30"
$ws.Range("N147").Value = 3.59
$ws.Range("O147").Value = "26/11/2023 18:13"
$ws.Range("P147").Value = 3.33
$ws.Range("Q147").Value = "02/12/2023 19:30"
$ws.Range("R147").Value = 1.95
$ws.Range("S147").Value = "26/11/2023 18:13"
$ws.Range("T147").Value = 2.36
$ws.Range("U147").Value = "02/12/2023 19:30"
$ws.Range("V147").Value = "https://www.betexplorer.com/football/poland/ekstraklasa/korona-kielce-lech-poznan/lEhvoSrr/"

$excel.CutCopyMode = $false
